$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (total) sheet and
#    populate it with the quarter's fund-holding detail rows.
#
#    NOTE: sheet handles in this host resolve by live position, not stable
#    identity, so once Worksheets.Add(before) shifts everything after the
#    insertion point, any handle obtained *before* the Add() (like the
#    "总计" handle used as the Before= argument) now tracks whatever sheet
#    sits at that old position - i.e. the newly inserted sheet. Always
#    re-fetch a handle by name *after* the sheet collection is mutated.
# ---------------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Column B (fund code) and columns D:G (scale / position figures) hold
# values that look numeric but must stay TEXT (leading zeros in fund codes
# such as "002624", decimal strings such as "12.98"). Format those ranges
# as Text before writing, otherwise Excel auto-converts the literals to
# numbers.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002624"
$newSheet.Range("C2").Value = "广发优企精选灵活配置混合A"
$newSheet.Range("D2").Value = "12.98"
$newSheet.Range("E2").Value = "92.40"
$newSheet.Range("F2").Value = "5.11"
$newSheet.Range("G2").Value = "0.6633"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "270025"
$newSheet.Range("C3").Value = "广发行业领先混合A"
$newSheet.Range("D3").Value = "11.11"
$newSheet.Range("E3").Value = "91.67"
$newSheet.Range("F3").Value = "5.31"
$newSheet.Range("G3").Value = "0.5899"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "960001"
$newSheet.Range("C4").Value = "广发行业领先混合H"
$newSheet.Range("D4").Value = "11.11"
$newSheet.Range("E4").Value = "91.67"
$newSheet.Range("F4").Value = "5.31"
$newSheet.Range("G4").Value = "0.5899"
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "501070"
$newSheet.Range("C5").Value = "广发睿阳三年定期开放混合"
$newSheet.Range("D5").Value = "7.06"
$newSheet.Range("E5").Value = "50.14"
$newSheet.Range("F5").Value = "6.01"
$newSheet.Range("G5").Value = "0.4243"
$newSheet.Range("H5").Value = 2

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "270008"
$newSheet.Range("C6").Value = "广发核心精选混合"
$newSheet.Range("D6").Value = "8.08"
$newSheet.Range("E6").Value = "93.75"
$newSheet.Range("F6").Value = "3.86"
$newSheet.Range("G6").Value = "0.3119"
$newSheet.Range("H6").Value = 10

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "377150"
$newSheet.Range("C7").Value = "上投摩根健康品质生活混合"
$newSheet.Range("D7").Value = "4.25"
$newSheet.Range("E7").Value = "87.59"
$newSheet.Range("F7").Value = "3.63"
$newSheet.Range("G7").Value = "0.1543"
$newSheet.Range("H7").Value = 6

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "010617"
$newSheet.Range("C8").Value = "兴业消费精选混合A"
$newSheet.Range("D8").Value = "2.96"
$newSheet.Range("E8").Value = "68.95"
$newSheet.Range("F8").Value = "3.99"
$newSheet.Range("G8").Value = "0.1181"
$newSheet.Range("H8").Value = 7

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "000747"
$newSheet.Range("C9").Value = "广发逆向策略灵活配置混合"
$newSheet.Range("D9").Value = "1.25"
$newSheet.Range("E9").Value = "89.61"
$newSheet.Range("F9").Value = "5.05"
$newSheet.Range("G9").Value = "0.0631"
$newSheet.Range("H9").Value = 8

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "010618"
$newSheet.Range("C10").Value = "兴业消费精选混合C"
$newSheet.Range("D10").Value = "1.47"
$newSheet.Range("E10").Value = "68.95"
$newSheet.Range("F10").Value = "3.99"
$newSheet.Range("G10").Value = "0.0587"
$newSheet.Range("H10").Value = 7

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "010257"
$newSheet.Range("C11").Value = "天弘多利一年定期开放混合"
$newSheet.Range("D11").Value = "2.64"
$newSheet.Range("E11").Value = "26.41"
$newSheet.Range("F11").Value = "0.79"
$newSheet.Range("G11").Value = "0.0209"
$newSheet.Range("H11").Value = 9

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "010021"
$newSheet.Range("C12").Value = "广发优企精选灵活配置混合C"
$newSheet.Range("D12").Value = "0.15"
$newSheet.Range("E12").Value = "92.40"
$newSheet.Range("F12").Value = "5.11"
$newSheet.Range("G12").Value = "0.0077"
$newSheet.Range("H12").Value = 10

# Style the header row (B1:H1) and the A-column index cells (A2:A12) the
# same way the rest of the workbook formats its header / index cells: copy
# the formatting already used for the header of the "2021-Q4" sheet.
$srcHeader = $wb.Worksheets.Item("2021-Q4")
$srcHeader.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcHeader.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Prepend a new "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down by one. Re-fetch the sheet handle by name since the
#    sheet collection changed above.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

# Give the about-to-be-created A5 the same formatting as the existing index
# cells (A2:A4) before any values are written into it.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B5").Value = $ws.Range("B4").Value2
$ws.Range("C5").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("D4").Value2

$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("D3").Value2

$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
